# Auto-generated Excel COM-interop script applying the Behemoth_Profits value updates.
# Each row change below corresponds to one <row> hunk in the source diff; values are
# plain numeric literals (no formulas in the source), so we just overwrite the H-N cells.

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 2
$ws.Range("H2").Value = 840.6
$ws.Range("J2").Value = 2
$ws.Range("L2").Value = 2
$ws.Range("N2").Value = -228
# Row 17
$ws.Range("H17").Value = 1772.6364
$ws.Range("J17").Value = 1772.6364
$ws.Range("L17").Value = 5317.9092
$ws.Range("N17").Value = -5653.9092
# Row 37
$ws.Range("H37").Value = 0
$ws.Range("I37").Value = 0
$ws.Range("K37").Value = 0
$ws.Range("M37").ClearContents()
# Row 62
$ws.Range("H62").Value = 5443.7
$ws.Range("I62").Value = 2610
$ws.Range("J62").Value = 6152.125
$ws.Range("K62").Value = 2610
$ws.Range("L62").Value = 6152.125
$ws.Range("M62").Value = -1986
$ws.Range("N62").Value = -7400.125
# Row 65
$ws.Range("H65").Value = 5443.7
$ws.Range("I65").Value = 2610
$ws.Range("J65").Value = 6152.125
$ws.Range("K65").Value = 13050
$ws.Range("L65").Value = 30760.625
$ws.Range("M65").Value = -9930
$ws.Range("N65").Value = -37000.625
# Row 103
$ws.Range("H103").Value = 5203
$ws.Range("J103").Value = 5255
$ws.Range("L103").Value = 15765
$ws.Range("N103").Value = -16937
# Row 116
$ws.Range("H116").Value = 6818.4546
$ws.Range("I116").Value = 6387.875
$ws.Range("K116").Value = 6387.875
$ws.Range("M116").Value = -2945.875
# Row 131
$ws.Range("H131").Value = 3369.1177
$ws.Range("I131").Value = 1527.6
$ws.Range("K131").Value = 4582.799999999999
$ws.Range("M131").Value = 457.2000000000007
# Row 135
$ws.Range("H135").Value = 1775
$ws.Range("I135").Value = 700
$ws.Range("K135").Value = 6300
$ws.Range("M135").Value = -3765
# Row 138
$ws.Range("H138").Value = 2198.782
$ws.Range("I138").Value = 764.8461
$ws.Range("J138").Value = 2915.75
$ws.Range("K138").Value = 2294.5383
$ws.Range("L138").Value = 8747.25
$ws.Range("M138").Value = 2845.4617
$ws.Range("N138").Value = -19027.25

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 8337558
$ws.Range("I32").Value = 10205357
$ws.Range("K32").Value = 10205357
$ws.Range("M32").Value = -10205070
# Row 74
$ws.Range("H74").Value = 13005326
$ws.Range("I74").Value = 17859704
$ws.Range("J74").Value = 1678445.1
$ws.Range("K74").Value = 17859704
$ws.Range("L74").Value = 1678445.1
$ws.Range("M74").Value = -17858830
$ws.Range("N74").Value = -1680193.1
# Row 77
$ws.Range("H77").Value = 13005326
$ws.Range("I77").Value = 17859704
$ws.Range("J77").Value = 1678445.1
$ws.Range("K77").Value = 89298520
$ws.Range("L77").Value = 8392225.5
$ws.Range("M77").Value = -89294152
$ws.Range("N77").Value = -8400961.5

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 134
$ws.Range("H134").Value = 360181.03
$ws.Range("J134").Value = 1669521.1
$ws.Range("L134").Value = 5008563.300000001
$ws.Range("N134").Value = -5013633.300000001

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 69
$ws.Range("H69").Value = 79385
$ws.Range("J69").Value = 99180
$ws.Range("L69").Value = 99180
$ws.Range("N69").Value = -100678
# Row 72
$ws.Range("H72").Value = 79385
$ws.Range("J72").Value = 99180
$ws.Range("L72").Value = 297540
$ws.Range("N72").Value = -305028
# Row 87
$ws.Range("H87").Value = 70233.375
$ws.Range("J87").Value = 55373.4
$ws.Range("L87").Value = 55373.4
$ws.Range("N87").Value = -57745.4
# Row 90
$ws.Range("H90").Value = 70233.375
$ws.Range("J90").Value = 55373.4
$ws.Range("L90").Value = 166120.2
$ws.Range("N90").Value = -177976.2
# Row 132
$ws.Range("H132").Value = 3627.4
$ws.Range("I132").Value = 3185.4614
$ws.Range("K132").Value = 9556.3842
$ws.Range("M132").Value = -7026.3842
# Row 140
$ws.Range("H140").Value = 72050
$ws.Range("J140").Value = 100000
$ws.Range("L140").Value = 100000
$ws.Range("N140").Value = -110360
# Row 141
$ws.Range("H141").Value = 346885.8
$ws.Range("J141").Value = 374428.66
$ws.Range("L141").Value = 374428.66
$ws.Range("N141").Value = -384788.66

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 2244.1
$ws.Range("I5").Value = 2244.1
$ws.Range("K5").Value = 6732.299999999999
$ws.Range("M5").Value = -6620.299999999999
# Row 114
$ws.Range("H114").Value = 1067
$ws.Range("I114").Value = 1256.6
$ws.Range("J114").Value = 830
$ws.Range("K114").Value = 3769.8
$ws.Range("L114").Value = 2490
$ws.Range("M114").Value = -515.7999999999997
$ws.Range("N114").Value = -8998
# Row 122
$ws.Range("H122").Value = 2224.1292
$ws.Range("I122").Value = 513.63635
$ws.Range("J122").Value = 3164.9
$ws.Range("K122").Value = 4622.72715
$ws.Range("L122").Value = 28484.1
$ws.Range("M122").Value = -2172.72715
$ws.Range("N122").Value = -33384.10000000001
# Row 131
$ws.Range("H131").Value = 7874.1704
$ws.Range("J131").Value = 8038.2827
$ws.Range("L131").Value = 24114.8481
$ws.Range("N131").Value = -34194.8481
# Row 135
$ws.Range("H135").Value = 2244.1
$ws.Range("I135").Value = 2244.1
$ws.Range("K135").Value = 20196.9
$ws.Range("M135").Value = -17661.9

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 97
$ws.Range("H97").Value = 1908
$ws.Range("I97").Value = 2123.077
$ws.Range("K97").Value = 2123.077
$ws.Range("M97").Value = -1627.077
# Row 122
$ws.Range("H122").Value = 1524.2
$ws.Range("I122").Value = 1579.4166
$ws.Range("K122").Value = 4738.2498
$ws.Range("M122").Value = -2288.2498

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 40
$ws.Range("H40").Value = 4752.7617
$ws.Range("J40").Value = 5333.3335
$ws.Range("L40").Value = 5333.3335
$ws.Range("N40").Value = -5605.3335
# Row 62
$ws.Range("H62").Value = 40000
$ws.Range("J62").Value = 40000
$ws.Range("L62").Value = 40000
$ws.Range("N62").Value = -41248
# Row 65
$ws.Range("H65").Value = 40000
$ws.Range("J65").Value = 40000
$ws.Range("L65").Value = 120000
$ws.Range("N65").Value = -126240
# Row 74
$ws.Range("H74").Value = 144000
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 144000
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 144000
$ws.Range("M74").ClearContents()
$ws.Range("N74").Value = -145996
# Row 77
$ws.Range("H77").Value = 144000
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 144000
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 432000
$ws.Range("M77").ClearContents()
$ws.Range("N77").Value = -441984
# Row 132
$ws.Range("H132").Value = 785674.0600000001
$ws.Range("I132").Value = 27036.857
$ws.Range("K132").Value = 81110.571
$ws.Range("M132").Value = -78580.571
# Row 136
$ws.Range("H136").Value = 48060.242
$ws.Range("I136").Value = 6419.273
$ws.Range("J136").Value = 131342.19
$ws.Range("K136").Value = 19257.819
$ws.Range("L136").Value = 394026.57
$ws.Range("M136").Value = -16707.819
$ws.Range("N136").Value = -399126.57

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 126
$ws.Range("H126").Value = 1286.619
$ws.Range("I126").Value = 1165.95
$ws.Range("K126").Value = 3497.85
$ws.Range("M126").Value = -1027.85
# Row 132
$ws.Range("H132").Value = 288466.03
$ws.Range("I132").Value = 2587.6875
$ws.Range("K132").Value = 7763.0625
$ws.Range("M132").Value = -5233.0625

Write-Output "Applied 37 row updates across ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets."
